$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final data for rows 2-19 (A: Player, B: Position, C: Team)
$data = @(
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Jeremy Sochan", "SF,PF", "San Antonio Spurs"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Dorian Finney-Smith", "SF,PF,C", "Brooklyn Nets"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Khris Middleton", "SF", "Milwaukee Bucks")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
